$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "User"
$ws.Range("B1").Value = "Link"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "IP"

# Update data row 2
$ws.Range("A2").Value = "Staszek"
$ws.Range("B2").Value = "https://manpages.ubuntu.com/manpages/trusty/pl/man1/mc.1.html"
$ws.Range("C2").Value = "2023-10-18 19:43:36"
$ws.Range("D2").Value = "172.30.0.1"

# Remove old row 3 entirely
$ws.Rows.Item(3).Delete()
